$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 199, shifting rows 199:218 down to 200:219
$ws.Rows.Item(199).Insert()

# Populate the new row 199 with the new data record
$ws.Cells.Item(199, 1).Value = 8
$ws.Cells.Item(199, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(199, 3).Value = "Coquimbo"
$ws.Cells.Item(199, 4).Value = 45223
$ws.Cells.Item(199, 4).NumberFormat = $ws.Cells.Item(200, 4).NumberFormat
$ws.Cells.Item(199, 5).Value = 4
$ws.Cells.Item(199, 6).Value = 100112052
$ws.Cells.Item(199, 7).Value = "Albahaca"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 800
$ws.Cells.Item(199, 11).Value = 3300
$ws.Cells.Item(199, 12).Value = 3500
$ws.Cells.Item(199, 13).Value = 3400
$ws.Cells.Item(199, 14).Value = "`$/paquete"
$ws.Cells.Item(199, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(199, 16).Value = 3400
$ws.Cells.Item(199, 17).Value = 1
$ws.Cells.Item(199, 18).Value = "Hortaliza"
